$wb = $excel.ActiveWorkbook

# --- Sheet 1: "metadata_included" ---
$ws1 = $wb.Worksheets.Item("metadata_included")

# Header row (row 3): rename the "series ..." labels to "type ..."
$ws1.Range("D3").Value = "...Other. types ….........."
$ws1.Range("E3").Value = "type               A"
$ws1.Range("F3").Value = "type B"
$ws1.Range("G3").Value = "type C"
$ws1.Range("H3").Value = "type D"

# Age-band labels
$ws1.Range("C4").Value = "< 66"
$ws1.Range("C5").Value = "> 65"

# Update the active selection to match the author's new selection
$ws1.Range("D3:H3").Select()

# --- Sheet 2: "no_metadata" ---
$ws2 = $wb.Worksheets.Item("no_metadata")

# Header row (row 1): rename the "series ..." labels to "type ..."
$ws2.Range("D1").Value = "...Other. types ….........."
$ws2.Range("E1").Value = "type               A"
$ws2.Range("F1").Value = "type B"
$ws2.Range("G1").Value = "type C"
$ws2.Range("H1").Value = "type D"

# Age-band labels
$ws2.Range("C2").Value = "< 66"
$ws2.Range("C3").Value = "> 65"

# Update the active selection to match the author's new selection
$ws2.Range("H9").Select()
